# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by the handback report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2).
# This timestamp is shared with the de-de sheet's "Correspond Handoff Datetime"
# for the same file (they recorded the same handoff-generation instant).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 13:06:44"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 13:06:39"
$wsZhCn.Range("K2").Value = "2016-08-23 13:07:09"

# de-de sheet: "Correspond Handoff Datetime" for row 2 mirrors the Overview value
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 13:06:44"

# de-de sheet: "Correspond Handback DateTime" for row 2
$wsDeDe.Range("K2").Value = "2016-08-23 13:07:20"
